$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2590.889
$ws.Range("I4").Value = 2590.889
$ws.Range("K4").Value = 2590.889
$ws.Range("M4").Value = -2476.889

$ws.Range("H5").Value = 107.666664
$ws.Range("I5").Value = 93.0625
$ws.Range("K5").Value = 93.0625
$ws.Range("M5").Value = 21.9375

$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -216
$ws.Range("N18").ClearContents()

$ws.Range("H32").Value = 8999.799999999999
$ws.Range("I32").Value = 4999.5
$ws.Range("K32").Value = 4999.5
$ws.Range("M32").Value = -4673.5

$ws.Range("H40").Value = 8347.799999999999
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 8347.799999999999
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 8347.799999999999
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -8697.799999999999

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

$ws.Range("H55").Value = 55
$ws.Range("J55").Value = 40
$ws.Range("L55").Value = 40
$ws.Range("N55").Value = -468

$ws.Range("H80").Value = 620
$ws.Range("I80").Value = 530
$ws.Range("J80").Value = 653.75
$ws.Range("K80").Value = 1590
$ws.Range("L80").Value = 1961.25
$ws.Range("M80").Value = -592
$ws.Range("N80").Value = -3957.25

$ws.Range("H83").Value = 620
$ws.Range("I83").Value = 530
$ws.Range("J83").Value = 653.75
$ws.Range("K83").Value = 4770
$ws.Range("L83").Value = 5883.75
$ws.Range("M83").Value = 222
$ws.Range("N83").Value = -15867.75

$ws.Range("H135").Value = 1137.3846
$ws.Range("I135").Value = 768.8182
$ws.Range("J135").Value = 3164.5
$ws.Range("K135").Value = 6919.3638
$ws.Range("L135").Value = 28480.5
$ws.Range("M135").Value = -4384.3638
$ws.Range("N135").Value = -33550.5

$ws.Range("H138").Value = 11499.5
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2864.0938
$ws.Range("I45").Value = 2605.5715
$ws.Range("K45").Value = 2605.5715
$ws.Range("M45").Value = -2228.5715

$ws.Range("H88").Value = 1516.6666
$ws.Range("J88").Value = 1525
$ws.Range("L88").Value = 1525
$ws.Range("N88").Value = -2337

$ws.Range("H91").Value = 1516.6666
$ws.Range("J91").Value = 1525
$ws.Range("L91").Value = 1525
$ws.Range("N91").Value = -4333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 616.3333
$ws.Range("I22").Value = 249.5
$ws.Range("J22").Value = 799.75
$ws.Range("K22").Value = 249.5
$ws.Range("L22").Value = 799.75
$ws.Range("M22").Value = -76.5
$ws.Range("N22").Value = -1145.75

$ws.Range("H86").Value = 1726.5454
$ws.Range("I86").Value = 2084.8572
$ws.Range("J86").Value = 1099.5
$ws.Range("K86").Value = 2084.8572
$ws.Range("L86").Value = 1099.5
$ws.Range("M86").Value = -961.8571999999999
$ws.Range("N86").Value = -3345.5

$ws.Range("H89").Value = 1726.5454
$ws.Range("I89").Value = 2084.8572
$ws.Range("J89").Value = 1099.5
$ws.Range("K89").Value = 10424.286
$ws.Range("L89").Value = 5497.5
$ws.Range("M89").Value = -4808.286
$ws.Range("N89").Value = -16729.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 18139
$ws.Range("I60").Value = 7847.5
$ws.Range("K60").Value = 7847.5
$ws.Range("M60").Value = -7336.5

$ws.Range("H62").Value = 2666.6667
$ws.Range("I62").Value = 2500
$ws.Range("J62").Value = 2750
$ws.Range("K62").Value = 2500
$ws.Range("L62").Value = 2750
$ws.Range("M62").Value = -1876
$ws.Range("N62").Value = -3998

$ws.Range("H65").Value = 2666.6667
$ws.Range("I65").Value = 2500
$ws.Range("J65").Value = 2750
$ws.Range("K65").Value = 12500
$ws.Range("L65").Value = 13750
$ws.Range("M65").Value = -9380
$ws.Range("N65").Value = -19990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 125237
$ws.Range("I4").Value = 251.4
$ws.Range("J4").Value = 333546.34
$ws.Range("K4").Value = 754.2
$ws.Range("L4").Value = 1000639.02
$ws.Range("M4").Value = -642.2
$ws.Range("N4").Value = -1000863.02

$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H81").Value = 2250
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 2250
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 6750
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -8996

$ws.Range("H84").Value = 2250
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 2250
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 20250
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -31482

$ws.Range("H136").Value = 4099.3335
$ws.Range("I136").Value = 3399.5
$ws.Range("J136").Value = 5499
$ws.Range("K136").Value = 10198.5
$ws.Range("L136").Value = 16497
$ws.Range("M136").Value = -5098.5
$ws.Range("N136").Value = -26697

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 4250
$ws.Range("I138").Value = 4250
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 12750
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -7610
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730

$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5140.6
$ws.Range("I7").Value = 4425.75
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 4425.75
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -4313.75
$ws.Range("N7").Value = -8224

$ws.Range("H16").Value = 350
$ws.Range("I16").Value = 350
$ws.Range("K16").Value = 350
$ws.Range("M16").Value = -180

$ws.Range("H22").Value = 900
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -405

$ws.Range("H27").Value = 900
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = -593

$ws.Range("H55").Value = 2655.111
$ws.Range("I55").Value = 4216.3335
$ws.Range("J55").Value = 1874.5
$ws.Range("K55").Value = 4216.3335
$ws.Range("L55").Value = 1874.5
$ws.Range("M55").Value = -4043.3335
$ws.Range("N55").Value = -2220.5

$ws.Range("H93").Value = 1964.6666
$ws.Range("I93").Value = 1964.6666
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1964.6666
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -716.6666
$ws.Range("N93").ClearContents()

$ws.Range("H126").Value = 5140.6
$ws.Range("I126").Value = 4425.75
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 13277.25
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -10807.25
$ws.Range("N126").Value = -28940

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31144

$ws.Range("H126").Value = 2180.6
$ws.Range("I126").Value = 2180.6
$ws.Range("K126").Value = 6541.799999999999
$ws.Range("M126").Value = -4071.799999999999

$ws.Range("H132").Value = 7885.222
$ws.Range("I132").Value = 4192.778
$ws.Range("J132").Value = 11577.667
$ws.Range("K132").Value = 12578.334
$ws.Range("L132").Value = 34733.001
$ws.Range("M132").Value = -10048.334
$ws.Range("N132").Value = -39793.001

$ws.Range("H135").Value = 167999.8
$ws.Range("J135").Value = 167999.8
$ws.Range("L135").Value = 167999.8
$ws.Range("N135").Value = -178139.8
